$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New "SIC" model tab (last sheet) gets populated with the impact-assessment
# indicator grid (Echelle / Critère / Indicateurs / ... table).
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SIC")

$data = @(
    @("Echelle", "Critère", "Indicateurs", "Valeur EI", "Justification prédiction court terme", "Incertitudes", "Valeur après impact/ MC CT", "Justification prédiction long terme", "Incertitudes", "Valeur après impact/ MC LT"),
    @("SITE", "Diversité espèce", "Nombre d'espèces dépendantes de l'habitat pour leur cycle de vie (Laisser la possibilité d’ajouter des taxons)"),
    @("SITE", "Diversité espèce", "Nombre d'espèces flore"),
    @("SITE", "Fonctionnalité", "Surface totale habitat"),
    @("SITE", "Fonctionnalité", "Nombre de patches d'habitat"),
    @("SITE", "Fonctionnalité", "Nombre de micro-habitats"),
    @("SITE", "Fonctionnalité", "Nombre d'horizons de sol par rapport à la référence"),
    @("SITE", "Fonctionnalité", "Epaisseur d'horizons organiques par rapport à la référence"),
    @("SITE", "Fonctionnalité", "Abondance relative de faune détritivore"),
    @("SITE", "Fonctionnalité", "Nombre d'espèces dépendantes de l'habitat pour la reproduction (Laisser la possibilité d’ajouter des taxons)"),
    @("SITE", "Fonctionnalité", "Nombre de Très Gros Bois Vivant"),
    @("SITE", "Fonctionnalité", "Proportion de bois mort"),
    @("SITE", "Fonctionnalité", "Nombre d'espèces bio-indicatrices"),
    @("SITE", "Fonctionnalité", "Densité de lichen"),
    @("SITE", "Fonctionnalité", "Ancienneté de la forêt"),
    @("SITE", "Fonctionnalité", "Nombre d'espèces pollinisatrices"),
    @("SITE", "Fonctionnalité", "Renvoyer vers la méthode sur les fonctionnalités de l'ONEMA"),
    @("SITE", "Structure", "Proportion de flore dominante"),
    @("SITE", "Structure", "Nombre de strates de végétation"),
    @("SITE", "Structure", "Hauteur de chaque strate"),
    @("SITE", "Pression", "Proportion de sol dégradé"),
    @("SITE", "Pression", "Nombre d'espèces indicatrices de pression"),
    @("SITE", "Pression", "Temps depuis la dernière coupe"),
    @("SITE", "Pression", "Taux de recouvrement des ligneux"),
    @("SITE", "Pression", "Taux de couvert des algues dues à l'eutrophisation"),
    @("ELARGI", "Connectivité", "Indice de fragmentation du type d'habitat"),
    @("ELARGI", "Représentativité", "Surface d'habitat dans le PE"),
    @("ELARGI", "Représentativité", "% surfacique d'habitat dans le PE")

)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $i + 1
    for ($j = 0; $j -lt $row.Count; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $row[$j]
    }
}

# Column widths to match the authored layout (characters; chosen so the
# stored OOXML width lands on - or as near as the 1/6-char pixel grid
# allows to - the target bestFit widths B:J = 15.86/99.71/19/32.71/11.71/
# 25.57/32/11.43/25.29).
$ws.Columns.Item(2).ColumnWidth = 15.0
$ws.Columns.Item(3).ColumnWidth = 98.83333333333333
$ws.Columns.Item(4).ColumnWidth = 18.166666666666668
$ws.Columns.Item(5).ColumnWidth = 31.833333333333332
$ws.Columns.Item(6).ColumnWidth = 10.833333333333334
$ws.Columns.Item(7).ColumnWidth = 24.666666666666668
$ws.Columns.Item(8).ColumnWidth = 31.166666666666668
$ws.Columns.Item(9).ColumnWidth = 10.666666666666666
$ws.Columns.Item(10).ColumnWidth = 24.5

# Selection on this sheet covers D2:D28, and it becomes the active/selected tab.
$ws.Range("D2:D28").Select()
$ws.Activate()
